$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert the first block of 4 new rows (before the old row 52 / "water_drill" row) ---
$ws.Range("A52:A55").EntireRow.Insert()

$ws.Range("A52").Value = "gui/hud/building_description/tower_cannon_acid"
$ws.Range("B52").Value = "Fires high caliber rounds dealing significant damage per shot. Acid coatied rounds deal damage over time"

$ws.Range("A53").Value = "gui/hud/building_description/tower_cannon_cooled"
$ws.Range("B53").Value = "Fires high caliber rounds dealing significant damage per shot. Cooling enables higher fire rate"

$ws.Range("A54").Value = "gui/hud/building_description/tower_cannon_cryo"
$ws.Range("B54").Value = "Fires high caliber rounds dealing significant damage per shot. Cryo infused rounds slow targets"

$ws.Range("A55").Value = "gui/hud/building_description/tower_cannon_incindiary"
$ws.Range("B55").Value = "Fires high caliber rounds dealing significant damage per shot. Plasma charged rounds deal fire damage"

# --- Insert the second block of 4 new rows (before the old row 98, now shifted to row 102) ---
$ws.Range("A102:A105").EntireRow.Insert()

$ws.Range("A102").Value = "gui/hud/building_name/tower_cannon_acid"
$ws.Range("B102").Value = "90mm Gun Tower, Acidic Rounds"

$ws.Range("A103").Value = "gui/hud/building_name/tower_cannon_cooled"
$ws.Range("B103").Value = "90mm Gun Tower, Cooled"

$ws.Range("A104").Value = "gui/hud/building_name/tower_cannon_cryo"
$ws.Range("B104").Value = "90mm Gun Tower, Cryo Rounds"

$ws.Range("A105").Value = "gui/hud/building_name/tower_cannon_incindiary"
$ws.Range("B105").Value = "90mm Gun Tower, Incidiary Rounds"

# --- Cosmetic view changes: zoom, frozen-pane scroll position, and final selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("A105").Select()
